# Commit: New crime data collected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings in row 8 / row 9) ---
# "Volume 32   Number  10" -> "...11"
$ws.Range("A8").Characters(21, 2).Text = "11"
# "Report Covering the Week  3/3/2025  Through  3/9/2025"
# -> "...3/10/2025  Through  3/16/2025" (edit 2nd date first so the 1st date offset does not shift)
$ws.Range("C9").Characters(46, 8).Text = "3/16/2025"
$ws.Range("C9").Characters(27, 8).Text = "3/10/2025"

# --- Cells whose type flips number -> text placeholder: copy format+text from a donor cell ---
# that already holds the matching placeholder string ("0" or "***.*"), then the row-specific
# numeric neighbours are overwritten below in the normal value pass.
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))

# --- Cells whose type flips text placeholder -> number: copy numeric format from a donor cell
# in the same row/column family, then overwrite with the real value below. ---
$ws.Range("F28").Copy($ws.Range("C28"))
$ws.Range("F28").Copy($ws.Range("D28"))
$ws.Range("H28").Copy($ws.Range("E28"))
$ws.Range("J29").Copy($ws.Range("D29"))
$ws.Range("J29").Copy($ws.Range("G29"))
$ws.Range("K29").Copy($ws.Range("E29"))
$ws.Range("K29").Copy($ws.Range("H29"))
$ws.Range("J30").Copy($ws.Range("D30"))
$ws.Range("J30").Copy($ws.Range("G30"))
$ws.Range("K30").Copy($ws.Range("E30"))
$ws.Range("K30").Copy($ws.Range("H30"))

# --- Plain value updates (no type/style change) ---
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -40
$ws.Range("C16").Value = 5
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 21
$ws.Range("K16").Value = 75
$ws.Range("L16").Value = 40
$ws.Range("M16").Value = -73.076923076923
$ws.Range("N16").Value = -90
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 37
$ws.Range("J17").Value = 64
$ws.Range("K17").Value = -42.1875
$ws.Range("L17").Value = -13.953488372093
$ws.Range("M17").Value = -32.727272727272
$ws.Range("N17").Value = -47.142857142857
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 22.222222222222
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 28
$ws.Range("L18").Value = 6.666666666666
$ws.Range("M18").Value = -56.164383561643
$ws.Range("N18").Value = -87.878787878787
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = -21.212121212121
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 83
$ws.Range("K19").Value = -12.048192771084
$ws.Range("L19").Value = -5.194805194805
$ws.Range("M19").Value = -15.116279069767
$ws.Range("N19").Value = -27.722772277227
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 10
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = -24.390243902439
$ws.Range("L20").Value = -3.125
$ws.Range("M20").Value = -60.25641025641
$ws.Range("N20").Value = -95.324283559577
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 23.529411764705
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -7.894736842105
$ws.Range("I21").Value = 198
$ws.Range("J21").Value = 230
$ws.Range("K21").Value = -13.91304347826
$ws.Range("L21").Value = -0.502512562814
$ws.Range("M21").Value = -47.757255936675
$ws.Range("N21").Value = -84.965831435079
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = -38.461538461538
$ws.Range("I24").Value = 119
$ws.Range("J24").Value = 155
$ws.Range("K24").Value = -23.225806451612
$ws.Range("L24").Value = -33.888888888888
$ws.Range("M24").Value = -22.222222222222
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = -22.222222222222
$ws.Range("I25").Value = 27
$ws.Range("J25").Value = 37
$ws.Range("K25").Value = -27.027027027027
$ws.Range("L25").Value = -15.625
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -26.666666666666
$ws.Range("I26").Value = 82
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 17.142857142857
$ws.Range("L26").Value = 26.153846153846
$ws.Range("M26").Value = -38.805970149253
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -62.5
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = -10
$ws.Range("L28").Value = 200
$ws.Range("J29").Value = 2
$ws.Range("J30").Value = 2

# --- Final values for the cells that changed type (set AFTER the format-donor Copy above) ---
$ws.Range("D16").Value = "0"
$ws.Range("E16").Value = "***.*"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
